$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C = IPC PO (neural network predicted output),
# Column D = DELTA (C - B), Column E = DELTA^2
$ws.Cells.Item(2, 3).Value = 30.36665211731184
$ws.Cells.Item(2, 4).Value = 0.8166521173118433
$ws.Cells.Item(2, 5).Value = 0.6669206807099166
$ws.Cells.Item(3, 3).Value = 30.37306242910338
$ws.Cells.Item(3, 4).Value = 0.6230624291033848
$ws.Cells.Item(3, 5).Value = 0.3882067905602105
$ws.Cells.Item(4, 3).Value = 29.99066224781861
$ws.Cells.Item(4, 4).Value = 0.1506622478186053
$ws.Cells.Item(4, 5).Value = 0.02269911291775483
$ws.Cells.Item(5, 3).Value = 29.91958244196424
$ws.Cells.Item(5, 4).Value = 0.1095824419642391
$ws.Cells.Item(5, 5).Value = 0.01200831158684582
$ws.Cells.Item(6, 3).Value = 29.89798823353122
$ws.Cells.Item(6, 4).Value = -0.02201176646878622
$ws.Cells.Item(6, 5).Value = 0.0004845178630763814
$ws.Cells.Item(7, 3).Value = 30.09988692223648
$ws.Cells.Item(7, 4).Value = 0.1198869222364785
$ws.Cells.Item(7, 5).Value = 0.01437287412333543
$ws.Cells.Item(8, 3).Value = 30.16321211180252
$ws.Cells.Item(8, 4).Value = 0.1232121118025233
$ws.Cells.Item(8, 5).Value = 0.01518122449483751
$ws.Cells.Item(9, 3).Value = 30.29750673130906
$ws.Cells.Item(9, 4).Value = 0.0875067313090625
$ws.Cells.Item(9, 5).Value = 0.00765742802439646
$ws.Cells.Item(10, 3).Value = 30.44672696430774
$ws.Cells.Item(10, 4).Value = 0.226726964307737
$ws.Cells.Item(10, 5).Value = 0.05140511634420183
$ws.Cells.Item(11, 3).Value = 30.29502784300016
$ws.Cells.Item(11, 4).Value = -0.08497215699984295
$ws.Cells.Item(11, 5).Value = 0.00722026746520596
$ws.Cells.Item(12, 3).Value = 30.3556783217689
$ws.Cells.Item(12, 4).Value = -0.08432167823110248
$ws.Cells.Item(12, 5).Value = 0.007110145419709582
$ws.Cells.Item(13, 3).Value = 30.53438401314706
$ws.Cells.Item(13, 4).Value = 0.05438401314706098
$ws.Cells.Item(13, 5).Value = 0.002957620885979702
$ws.Cells.Item(14, 3).Value = 30.63601958917854
$ws.Cells.Item(14, 4).Value = -0.05398041082146321
$ws.Cells.Item(14, 5).Value = 0.002913884752453942
$ws.Cells.Item(15, 3).Value = 30.17980980881382
$ws.Cells.Item(15, 4).Value = -0.5701901911861782
$ws.Cells.Item(15, 5).Value = 0.3251168541249305
$ws.Cells.Item(16, 3).Value = 30.30393750451478
$ws.Cells.Item(16, 4).Value = -0.6360624954852234
$ws.Cells.Item(16, 5).Value = 0.4045754981628898
$ws.Cells.Item(17, 3).Value = 30.58400968609432
$ws.Cells.Item(17, 4).Value = -0.3659903139056802
$ws.Cells.Item(17, 5).Value = 0.1339489098727783
$ws.Cells.Item(18, 3).Value = 30.87614945704727
$ws.Cells.Item(18, 4).Value = -0.1438505429527304
$ws.Cells.Item(18, 5).Value = 0.02069297870779532
$ws.Cells.Item(19, 3).Value = 31.1134232730596
$ws.Cells.Item(19, 4).Value = -0.006576726940398459
$ws.Cells.Item(19, 5).Value = 0.000043253337248562871339939429
$ws.Cells.Item(20, 3).Value = 31.16860810802609
$ws.Cells.Item(20, 4).Value = -0.1113918919739092
$ws.Cells.Item(20, 5).Value = 0.01240815359752706
$ws.Cells.Item(21, 3).Value = 31.33760770974518
$ws.Cells.Item(21, 4).Value = -0.04239229025482061
$ws.Cells.Item(21, 5).Value = 0.001797106273048958
$ws.Cells.Item(22, 3).Value = 31.63689533939427
$ws.Cells.Item(22, 4).Value = 0.05689533939427349
$ws.Cells.Item(22, 5).Value = 0.003237079644789569
$ws.Cells.Item(23, 3).Value = 31.57125835059635
$ws.Cells.Item(23, 4).Value = -0.07874164940365347
$ws.Cells.Item(23, 5).Value = 0.006200247350807881
$ws.Cells.Item(24, 3).Value = 31.902911412466
$ws.Cells.Item(24, 4).Value = 0.02291141246600503
$ws.Cells.Item(24, 5).Value = 0.0005249328211874108
$ws.Cells.Item(25, 3).Value = 32.27999587101024
$ws.Cells.Item(25, 4).Value = -0.000004128989758100942708551884
$ws.Cells.Item(25, 5).Value = 0.000000000017048556422502481384
$ws.Cells.Item(26, 3).Value = 32.31492129653395
$ws.Cells.Item(26, 4).Value = -0.135078703466057
$ws.Cells.Item(26, 5).Value = 0.01824625613007096
$ws.Cells.Item(27, 3).Value = 33.33681736660397
$ws.Cells.Item(27, 4).Value = 0.4868173666039723
$ws.Cells.Item(27, 5).Value = 0.2369911484272264
$ws.Cells.Item(28, 3).Value = 33.02364523339497
$ws.Cells.Item(28, 4).Value = 0.1236452333949742
$ws.Cells.Item(28, 5).Value = 0.01528814374129765
$ws.Cells.Item(29, 3).Value = 33.37918762896857
$ws.Cells.Item(29, 4).Value = 0.2791876289685646
$ws.Cells.Item(29, 5).Value = 0.07794573216908891
$ws.Cells.Item(30, 3).Value = 33.43665432956772
$ws.Cells.Item(30, 4).Value = 0.03665432956772463
$ws.Cells.Item(30, 5).Value = 0.001343539876059372
$ws.Cells.Item(31, 3).Value = 33.98856141008333
$ws.Cells.Item(31, 4).Value = 0.2885614100833251
$ws.Cells.Item(31, 5).Value = 0.08326768738927694
$ws.Cells.Item(32, 3).Value = 34.70679810123244
$ws.Cells.Item(32, 4).Value = 0.6067981012324353
$ws.Cells.Item(32, 5).Value = 0.3682039356592888
$ws.Cells.Item(33, 3).Value = 34.64649225017088
$ws.Cells.Item(33, 4).Value = 0.2464922501708813
$ws.Cells.Item(33, 5).Value = 0.06075842939430432
$ws.Cells.Item(34, 3).Value = 35.03841452557509
$ws.Cells.Item(34, 4).Value = 0.1384145255750937
$ws.Cells.Item(34, 5).Value = 0.01915858089017827
$ws.Cells.Item(35, 3).Value = 34.9147971642647
$ws.Cells.Item(35, 4).Value = -0.3852028357352992
$ws.Cells.Item(35, 5).Value = 0.1483812246585159
$ws.Cells.Item(36, 3).Value = 35.13928856504177
$ws.Cells.Item(36, 4).Value = -0.560711434958229
$ws.Cells.Item(36, 5).Value = 0.3143973132929162
$ws.Cells.Item(37, 3).Value = 35.700520038417
$ws.Cells.Item(37, 4).Value = -0.5994799615829933
$ws.Cells.Item(37, 5).Value = 0.3593762243395471
$ws.Cells.Item(38, 3).Value = 36.07152366051772
$ws.Cells.Item(38, 4).Value = -0.7284763394822775
$ws.Cells.Item(38, 5).Value = 0.5306777771854985
$ws.Cells.Item(39, 3).Value = 37.06429871258559
$ws.Cells.Item(39, 4).Value = -0.2357012874144075
$ws.Cells.Item(39, 5).Value = 0.05555509688880915
$ws.Cells.Item(40, 3).Value = 37.59373788371554
$ws.Cells.Item(40, 4).Value = -0.3062621162844579
$ws.Cells.Item(40, 5).Value = 0.0937964838710348
$ws.Cells.Item(41, 3).Value = 38.34942907022938
$ws.Cells.Item(41, 4).Value = -0.1505709297706233
$ws.Cells.Item(41, 5).Value = 0.02267160489198996
$ws.Cells.Item(42, 3).Value = 39.32514080133006
$ws.Cells.Item(42, 4).Value = 0.425140801330059
$ws.Cells.Item(42, 5).Value = 0.1807447009555647
$ws.Cells.Item(43, 3).Value = 40.0991865390098
$ws.Cells.Item(43, 4).Value = 0.6991865390098013
$ws.Cells.Item(43, 5).Value = 0.4888618163325044
$ws.Cells.Item(44, 3).Value = 40.5984031086222
$ws.Cells.Item(44, 4).Value = 0.6984031086222018
$ws.Cells.Item(44, 5).Value = 0.487766902133155
$ws.Cells.Item(45, 3).Value = 40.20863638200954
$ws.Cells.Item(45, 4).Value = 0.1086363820095357
$ws.Cells.Item(45, 5).Value = 0.01180186349612178
$ws.Cells.Item(46, 3).Value = 41.23139036539403
$ws.Cells.Item(46, 4).Value = 0.6313903653940329
$ws.Cells.Item(46, 5).Value = 0.3986537935124104
$ws.Cells.Item(47, 3).Value = 41.46436510051771
$ws.Cells.Item(47, 4).Value = 0.5643651005177119
$ws.Cells.Item(47, 5).Value = 0.3185079666823671
$ws.Cells.Item(48, 3).Value = 41.32579206646838
$ws.Cells.Item(48, 4).Value = 0.125792066468378
$ws.Cells.Item(48, 5).Value = 0.01582364398638484
$ws.Cells.Item(49, 3).Value = 40.99730796859584
$ws.Cells.Item(49, 4).Value = -0.502692031404159
$ws.Cells.Item(49, 5).Value = 0.2526992784372399
$ws.Cells.Item(50, 3).Value = 41.2627089484069
$ws.Cells.Item(50, 4).Value = -0.537291051593094
$ws.Cells.Item(50, 5).Value = 0.2886816741220128
$ws.Cells.Item(51, 3).Value = 41.65376525386564
$ws.Cells.Item(51, 4).Value = -0.5462347461343668
$ws.Cells.Item(51, 5).Value = 0.2983723978844762

# TOTAL row (52): C = sum of DELTA, E = sum of DELTA^2
$ws.Cells.Item(52, 3).Value = 0.9627802583703939
$ws.Cells.Item(52, 5).Value = 7.255656205405317

# MSE row (53): E = mean of DELTA^2
$ws.Cells.Item(53, 5).Value = 0.1451131241081063
